# Insert a new data row at row 609 (pushing the existing rows 609..729 down
# to 610..730) and populate the newly inserted row with the new record.
# This matches the diff: dimension grows from A1:T729 to A1:T730, and every
# previously-existing row at or after 609 shifts down by one row, while the
# brand new row 609 carries the new "Primera" / Perú entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 609 downward by inserting a fresh row above the current row 609.
$ws.Rows.Item(609).Insert()

# Populate the newly inserted row 609 with the new record's values.
$ws.Cells.Item(609, 1).Value = 6
$ws.Cells.Item(609, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(609, 3).Value = "Metropolitana"
$ws.Cells.Item(609, 4).Value = 45173
$ws.Cells.Item(609, 5).Value = 13
$ws.Cells.Item(609, 6).Value = "Fruta"
$ws.Cells.Item(609, 7).Value = 100101
$ws.Cells.Item(609, 8).Value = "Berries"
$ws.Cells.Item(609, 9).Value = 100101001
$ws.Cells.Item(609, 10).Value = "Arándano (blue)"
$ws.Cells.Item(609, 11).Value = "Sin especificar"
$ws.Cells.Item(609, 12).Value = "Primera"
$ws.Cells.Item(609, 13).Value = 365
$ws.Cells.Item(609, 14).Value = 12000
$ws.Cells.Item(609, 15).Value = 12000
$ws.Cells.Item(609, 16).Value = 12000
$ws.Cells.Item(609, 17).Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Cells.Item(609, 18).Value = "Perú"
$ws.Cells.Item(609, 19).Value = 8000
$ws.Cells.Item(609, 20).Value = 1.5
